$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
Set-TextValue $ws "D2" "65.381.79"
Set-TextValue $ws "E2" "  +0.02%  "

# Row 3
Set-TextValue $ws "D3" "3.545.19"
Set-TextValue $ws "E3" "  +3.43%  "

# Row 4
Set-TextValue $ws "D4" "0.999"
Set-TextValue $ws "E4" "  -0.09%  "

# Row 5
Set-TextValue $ws "D5" "597.09"
Set-TextValue $ws "E5" "  +0.45%  "

# Row 6
Set-TextValue $ws "D6" "139.82"
Set-TextValue $ws "E6" "  +4.39%  "

# Row 7
Set-TextValue $ws "D7" "3.546.58"
Set-TextValue $ws "E7" "  +3.58%  "

# Row 8
Set-TextValue $ws "E8" "  +0.10%  "

# Row 9
Set-TextValue $ws "E9" "  +1.42%  "

# Row 10
Set-TextValue $ws "E10" "  +3.56%  "

# Row 11
Set-TextValue $ws "D11" "7.17"
Set-TextValue $ws "E11" "  -4.09%  "

# Row 12
Set-TextValue $ws "E12" "  +3.97%  "

# Row 13
Set-TextValue $ws "D13" "4.139.67"
Set-TextValue $ws "E13" "  +3.41%  "

# Row 14
Set-TextValue $ws "E14" "  +5.08%  "

# Row 15
Set-TextValue $ws "D15" "26.97"
Set-TextValue $ws "E15" "  +2.37%  "

# Row 16
Set-TextValue $ws "D16" "3.533.23"
Set-TextValue $ws "E16" "  +4.20%  "

# Row 17
Set-TextValue $ws "E17" "  +1.50%  "

# Row 18
Set-TextValue $ws "D18" "65.231.21"
Set-TextValue $ws "E18" "  -0.14%  "

# Row 19
Set-TextValue $ws "D19" "10.36"
Set-TextValue $ws "E19" "  +5.19%  "

# Row 20
Set-TextValue $ws "E20" "  +2.47%  "

# Row 21
Set-TextValue $ws "D21" "14.25"
Set-TextValue $ws "E21" "  +4.31%  "

# Row 22
Set-TextValue $ws "D22" "396.74"
Set-TextValue $ws "E22" "  +1.46%  "

# Row 23
Set-TextValue $ws "D23" "0.572"
Set-TextValue $ws "E23" "  +5.50%  "

# Row 24
Set-TextValue $ws "E24" "  +1.94%  "

# Row 25
Set-TextValue $ws "D25" "3.679.16"
Set-TextValue $ws "E25" "  +3.12%  "

# Row 26
Set-TextValue $ws "E26" "  +0.05%  "

# Row 27
Set-TextValue $ws "E27" "  +9.59%  "

# Row 28
Set-TextValue $ws "D28" "7.79"
Set-TextValue $ws "E28" "  +9.22%  "

# Row 29
Set-TextValue $ws "D29" "0.999"
Set-TextValue $ws "E29" "  -0.19%  "

# Row 30
Set-TextValue $ws "E30" "  +1.09%  "

# Row 31
Set-TextValue $ws "D31" "8.31"
Set-TextValue $ws "E31" "  +2.00%  "

# Row 32
Set-TextValue $ws "D32" "3.554.17"
Set-TextValue $ws "E32" "  +3.57%  "

# Row 33
Set-TextValue $ws "B33" "EthereumClassic"
Set-TextValue $ws "C33" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D33" "24.02"
Set-TextValue $ws "E33" "  +6.13%  "

# Row 34
Set-TextValue $ws "B34" "USDe"
Set-TextValue $ws "C34" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws "D34" "1.00"
Set-TextValue $ws "E34" "  +0.02%  "

# Row 35
Set-TextValue $ws "E35" "  +0.39%  "

# Row 36
Set-TextValue $ws "D36" "1.24"
Set-TextValue $ws "E36" "  +2.07%  "

# Row 37
Set-TextValue $ws "D37" "7.05"
Set-TextValue $ws "E37" "  +3.80%  "

# Row 38
Set-TextValue $ws "D38" "169.30"
Set-TextValue $ws "E38" "  -1.90%  "

# Row 39
Set-TextValue $ws "D39" "1.55"
Set-TextValue $ws "E39" "  +2.90%  "

# Row 40
Set-TextValue $ws "D40" "4.95"
Set-TextValue $ws "E40" "  +2.87%  "

# Row 41
Set-TextValue $ws "D41" "0.0807"
Set-TextValue $ws "E41" "  +5.12%  "

# Row 42
Set-TextValue $ws "D42" "0.825"
Set-TextValue $ws "E42" "  +1.81%  "

# Row 43
Set-TextValue $ws "D43" "26.73"
Set-TextValue $ws "E43" "  +21.92%  "

# Row 44
Set-TextValue $ws "D44" "42.70"
Set-TextValue $ws "E44" "  -1.63%  "

# Row 45
Set-TextValue $ws "E45" "  -0.05%  "

# Row 46
Set-TextValue $ws "E46" "  +1.70%  "

# Row 47
Set-TextValue $ws "D47" "1.20"
Set-TextValue $ws "E47" "  +10.66%  "

# Row 48
Set-TextValue $ws "E48" "  +4.12%  "

# Row 49
Set-TextValue $ws "D49" "6.85"
Set-TextValue $ws "E49" "  +5.32%  "

# Row 50
Set-TextValue $ws "D50" "2.380.26"
Set-TextValue $ws "E50" "  +8.97%  "

# Row 51
Set-TextValue $ws "B51" "dogwifhat"
Set-TextValue $ws "C51" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws "D51" "2.13"
Set-TextValue $ws "E51" "  +0.11%  "
